# Atualizei dados da bibi e add
# Insere um novo registro (dia 11 de agosto/2025) na planilha de
# faturamento diario, logo apos o ultimo dia de agosto existente (linha 11),
# deslocando os registros de julho/junho/maio uma linha para baixo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova linha na posicao 12, empurrando tudo que estava
# a partir dali (julho, junho, maio) uma linha para baixo.
$ws.Rows.Item(12).Insert(-4121)   # -4121 = xlShiftDown

# Preenche a nova linha 12 com o novo dia de agosto/2025
$ws.Cells.Item(12, 1).Value = 11
$ws.Cells.Item(12, 2).Value = 32841.8
$ws.Cells.Item(12, 3).Value = 8
$ws.Cells.Item(12, 4).Value = 2025
$ws.Cells.Item(12, 5).Value = "08/2025"
